$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 53
$ws.Cells.Item($row, 1).Value = "T come tigro"
$ws.Cells.Item($row, 2).Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Cells.Item($row, 3).Value = "Leonardo Viola | Shark Attack"
$ws.Cells.Item($row, 4).Value = "Alessio Bragagna | FC Savignano"
$ws.Cells.Item($row, 5).Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Cells.Item($row, 6).Value = "Roberto Barozzi | Demobusters"
